# Update TPM-derived metrics for the Cdh1-Egfr LR-pair worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = [double]"0.002166666666666667"
$ws.Range("H2").Value = [double]"0.0065"
$ws.Range("I2").Value = [double]"0.004890446475191893"
$ws.Range("J2").Value = [double]"0.004890446475191893"
$ws.Range("M2").Value = [double]"0.428743"
$ws.Range("N2").Value = [double]"1.286229"
$ws.Range("O2").Value = [double]"0.00412050394863168"
$ws.Range("P2").Value = [double]"0.00412050394863168"
$ws.Range("Q2").Value = [double]"0.0009289431666666667"
$ws.Range("R2").Value = [double]"0.008360488500000001"
$ws.Range("S2").Value = [double]"2.015110401160008E-05"
$ws.Range("T2").Value = [double]"2.015110401160008E-05"

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = [double]"0.002166666666666667"
$ws.Range("H3").Value = [double]"0.0065"
$ws.Range("I3").Value = [double]"0.004890446475191893"
$ws.Range("J3").Value = [double]"0.004890446475191893"
$ws.Range("M3").Value = [double]"80.22623699999998"
$ws.Range("O3").Value = [double]"0.7710272268990069"
$ws.Range("P3").Value = [double]"0.7710272268990069"
$ws.Range("Q3").Value = [double]"0.1738235134999999"
$ws.Range("R3").Value = [double]"1.5644116215"
$ws.Range("S3").Value = [double]"0.003770667384065229"
$ws.Range("T3").Value = [double]"0.003770667384065229"

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = [double]"0.002166666666666667"
$ws.Range("H4").Value = [double]"0.0065"
$ws.Range("I4").Value = [double]"0.004890446475191893"
$ws.Range("J4").Value = [double]"0.004890446475191893"
$ws.Range("M4").Value = [double]"23.39612766666667"
$ws.Range("N4").Value = [double]"70.188383"
$ws.Range("O4").Value = [double]"0.2248522691523614"
$ws.Range("P4").Value = [double]"0.2248522691523614"
$ws.Range("Q4").Value = [double]"0.05069160994444444"
$ws.Range("R4").Value = [double]"0.4562244895"
$ws.Range("S4").Value = [double]"0.001099627987115065"
$ws.Range("T4").Value = [double]"0.001099627987115065"

# Row 5 (MuSCs -> ECs)
$ws.Range("G5").Value = [double]"0.440874"
$ws.Range("H5").Value = [double]"1.322622"
$ws.Range("I5").Value = [double]"0.9951095535248081"
$ws.Range("J5").Value = [double]"0.9951095535248081"
$ws.Range("M5").Value = [double]"0.428743"
$ws.Range("N5").Value = [double]"1.286229"
$ws.Range("O5").Value = [double]"0.00412050394863168"
$ws.Range("P5").Value = [double]"0.00412050394863168"
$ws.Range("Q5").Value = [double]"0.189021641382"
$ws.Range("R5").Value = [double]"1.701194772438"
$ws.Range("S5").Value = [double]"0.00410035284462008"
$ws.Range("T5").Value = [double]"0.00410035284462008"

# Row 6 (MuSCs -> FAPs)
$ws.Range("G6").Value = [double]"0.440874"
$ws.Range("H6").Value = [double]"1.322622"
$ws.Range("I6").Value = [double]"0.9951095535248081"
$ws.Range("J6").Value = [double]"0.9951095535248081"
$ws.Range("M6").Value = [double]"80.22623699999998"
$ws.Range("O6").Value = [double]"0.7710272268990069"
$ws.Range("P6").Value = [double]"0.7710272268990069"
$ws.Range("Q6").Value = [double]"35.36966201113799"
$ws.Range("R6").Value = [double]"318.3269581002419"
$ws.Range("S6").Value = [double]"0.7672565595149418"
$ws.Range("T6").Value = [double]"0.7672565595149418"

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = [double]"0.440874"
$ws.Range("H7").Value = [double]"1.322622"
$ws.Range("I7").Value = [double]"0.9951095535248081"
$ws.Range("J7").Value = [double]"0.9951095535248081"
$ws.Range("M7").Value = [double]"23.39612766666667"
$ws.Range("N7").Value = [double]"70.188383"
$ws.Range("O7").Value = [double]"0.2248522691523614"
$ws.Range("P7").Value = [double]"0.2248522691523614"
$ws.Range("Q7").Value = [double]"10.314744388914"
$ws.Range("R7").Value = [double]"92.832699500226"
$ws.Range("S7").Value = [double]"0.2237526411652464"
$ws.Range("T7").Value = [double]"0.2237526411652464"
